$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the 'Posted At' column keeps storing plain text dates (not Excel date serials).
$ws.Range("F2:F17").NumberFormat = "@"

$jobs = @(
    ,@('Software Engineer III - AI/ML', 'JPMorganChase', 'Houston, TX, US USA', 21.1, 'Data Scientist, LangChain, RAG, LLaMA, Pinecone, Prompt Engineering, TensorFlow, PyTorch, Azure ML, Docker', '2026-02-26', 'https://www.indeed.com/viewjob?jk=330acb9fbf425d3a')
    ,@('Senior Data Engineer (Remote in CA)', 'First American', 'Santa Ana, CA, US USA', 20, 'Generative AI, RAG, Redshift, BigQuery, Data Lake, Docker, Kubernetes, Apache Airflow, Snowflake, Databricks', '2026-02-25', 'https://www.indeed.com/viewjob?jk=05ccb57abd2eb2be')
    ,@('DevOps & Backend Engineer — Mid/Senior', 'Open Insurance', 'Remote, US USA', 20, 'RAG, Pinecone, S3, EC2, Docker, Kubernetes, CI/CD, GitHub Actions, Terraform, Git', '2026-02-25', 'https://www.indeed.com/viewjob?jk=6691fd82745fd2c9')
    ,@('DevOps Software Engineer III', 'JPMorganChase', 'Plano, TX, US USA', 17.8, 'RAG, S3, EC2, Athena, Redshift, Docker, Kubernetes, Jenkins, Terraform, Redshift', '2026-02-26', 'https://www.indeed.com/viewjob?jk=ae9c515199307b36')
    ,@('Senior Machine Learning Engineer (Remote)', 'The Home Depot', 'Atlanta, GA, US USA', 17.8, 'Machine Learning Engineer, RAG, TensorFlow, PyTorch, Gensim, BigQuery, CI/CD, Git, BigQuery, NoSQL', '2026-02-26', 'https://www.indeed.com/viewjob?jk=64eeda5b2dfcfe99')
    ,@('Machine Learning Engineer II (Remote)', 'The Home Depot', 'Atlanta, GA, US USA', 16.7, 'Machine Learning Engineer, RAG, TensorFlow, PyTorch, Gensim, BigQuery, CI/CD, Git, BigQuery, Python', '2026-02-26', 'https://www.indeed.com/viewjob?jk=dd0ef78a415d3113')
    ,@('Software Engineer III - Data Engineering Java/Python', 'JPMorganChase', 'Jersey City, NJ, US USA', 14.4, 'RAG, Data Lake, CI/CD, Snowflake, Databricks, PySpark, Kafka, NoSQL, Python, SQL', '2026-02-26', 'https://www.indeed.com/viewjob?jk=5ec8999487afb4f6')
    ,@('Data Engineering Specialist (Hybrid)', 'Selective Insurance', 'Charlotte, NC, US USA', 12.2, 'Machine Learning Engineer, RAG, Synapse, Data Lake, Kinesis, CI/CD, Databricks, Kafka, Python, SQL', '2026-02-25', 'https://www.indeed.com/viewjob?jk=6bfccdbabce88c5e')
    ,@('Engineer II', 'Mutual of Omaha', 'Remote, US USA', 12.2, 'Docker, Kubernetes, CI/CD, Jenkins, GitHub Actions, Git, PostgreSQL, SQL, R, Java', '2026-02-26', 'https://www.indeed.com/viewjob?jk=e4eee0c9ff41c6e6')
    ,@('Backend & CyberSec Engineer — Mid/Senior', 'Open Insurance', 'Remote, US USA', 12.2, 'FastAPI, Docker, Kubernetes, CI/CD, Jenkins, GitHub Actions, Terraform, Git, Python, SQL', '2026-02-25', 'https://www.indeed.com/viewjob?jk=01c5cbaeb48f13ea')
    ,@('Machine Learning Engineer', 'Workday', 'Boulder, CO, US USA', 11.1, 'Machine Learning Engineer, RAG, TensorFlow, PyTorch, Docker, Kubernetes, PySpark, Python, R, Scala', '2026-02-26', 'https://www.indeed.com/viewjob?jk=cfc72a528e4a243f')
    ,@('Applied AIML Data Scientist - Associate', 'JPMorganChase', 'Jersey City, NJ, US USA', 10, 'Data Scientist, Machine Learning Engineer, RAG, TensorFlow, XGBoost, Keras, Python, SQL, R', '2026-02-26', 'https://www.indeed.com/viewjob?jk=530092c6eae96e74')
    ,@('Software Engineer III - Data Engineering and MLOps', 'JPMorganChase', 'Boston, MA, US USA', 10, 'RAG, Data Lake, CI/CD, Databricks, Python, SQL, R, Java, Scala', '2026-02-26', 'https://www.indeed.com/viewjob?jk=bf33b7d584ddd74d')
    ,@('Systematic Production Support Engineer', 'Millennium Management', 'Miami, FL, US USA', 10, 'Generative AI, Jenkins, Kafka, PostgreSQL, Python, SQL, R, Java, Scala', '2026-02-25', 'https://www.indeed.com/viewjob?jk=77dbe60d954cf27d')
    ,@('Senior Automation Enablement Engineer', 'Conagra Brands', 'Omaha, NE, US USA', 10, 'AI Engineer, Generative AI, Copilot, Prompt Engineering, CI/CD, Terraform, Python, R, Scala', '2026-02-25', 'https://www.indeed.com/viewjob?jk=1d83e0839c59ed13')
    ,@('Systematic Production Support Engineer', 'Millennium Management', 'New York, NY, US USA', 10, 'Generative AI, Jenkins, Kafka, PostgreSQL, Python, SQL, R, Java, Scala', '2026-02-20', 'https://www.indeed.com/viewjob?jk=d69e2a46fbb76bda')
)

$row = 2
foreach ($job in $jobs) {
    $ws.Cells.Item($row, 1).Value = $job[0]
    $ws.Cells.Item($row, 2).Value = $job[1]
    $ws.Cells.Item($row, 3).Value = $job[2]
    $ws.Cells.Item($row, 4).Value = $job[3]
    $ws.Cells.Item($row, 5).Value = $job[4]
    $ws.Cells.Item($row, 6).Value = $job[5]
    $ws.Cells.Item($row, 7).Value = $job[6]
    $row++
}
